$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10 have their image filename extensions changed from uppercase
# ".JPG" to lowercase ".jpg" in columns F, G and H (img1/img2/img3).
for ($r = 2; $r -le 10; $r++) {
    for ($c = 6; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $old = $cell.Value2
        $new = $old -replace "\.JPG$", ".jpg"
        $cell.Value = $new
    }
}

# Update the stored selection to match the final state of the workbook.
$ws.Range("F27").Select()
